# IKD update: GaN CMOS 2026-02-16T23:33Z
# Appends three new literature records (rows 197-199) to the "Master" sheet,
# matching the rows added in the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateTextCell($cell, $text) {
    # AddedDate values look like dates ("2026-02-16"). Force the cell to be
    # treated as literal text first so Excel does not auto-convert the
    # string into a date serial number (the source data stores these as
    # plain text, not as real dates).
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $text
}

# ---------------------------------------------------------------------------
# Row 197
# ---------------------------------------------------------------------------
$ws.Range("B197").Value = "Neuromorphic Multi-LLM Modular Intelligence Architecture: A Systems-Level Path Toward AGI Beyond Monolithic LLM Limits"
$ws.Range("C197").Value = 2026
$ws.Range("D197").Value = "Institute of Electrical and Electronics Engineers (IEEE)"
$ws.Range("F197").Value = "Mishra, Anindya"
$ws.Range("H197").Value = "10.36227/techrxiv.177127405.56163861/v1"
$ws.Range("I197").Value = "https://doi.org/10.36227/techrxiv.177127405.56163861/v1"
$ws.Range("J197").Value = "Journal"
$ws.Range("K197").Value = "Co-integration"
$ws.Range("L197").Value = "Experiment"
$ws.Range("M197").Value = "Contacts"
$ws.Range("Q197").Value = "Neuromorphic Multi-LLM Modular Intelligence Architecture: A Systems-Level Path Toward AGI Beyond Monolithic LLM Limits"
$ws.Range("R197").Value = "High"
Set-DateTextCell "S197" "2026-02-16"

# ---------------------------------------------------------------------------
# Row 198
# ---------------------------------------------------------------------------
$ws.Range("B198").Value = "Neuromorphic Multi-LLM Modular Intelligence Architecture: A Systems-Level Path Toward AGI Beyond Monolithic LLM Limits"
$ws.Range("C198").Value = 2026
$ws.Range("D198").Value = "Institute of Electrical and Electronics Engineers (IEEE)"
$ws.Range("F198").Value = "Mishra, Anindya"
$ws.Range("H198").Value = "10.36227/techrxiv.177127405.56163861/v1"
$ws.Range("I198").Value = "https://doi.org/10.36227/techrxiv.177127405.56163861/v1"
$ws.Range("J198").Value = "Journal"
$ws.Range("K198").Value = "Co-integration"
$ws.Range("L198").Value = "Experiment"
$ws.Range("M198").Value = "Contacts"
$ws.Range("Q198").Value = "Neuromorphic Multi-LLM Modular Intelligence Architecture: A Systems-Level Path Toward AGI Beyond Monolithic LLM Limits"
$ws.Range("R198").Value = "High"
Set-DateTextCell "S198" "2026-02-16"

# ---------------------------------------------------------------------------
# Row 199
# ---------------------------------------------------------------------------
$ws.Range("B199").Value = "Advances and Perspectives in Gate Dielectric Thin Films for 4H-SiC MOSFETs"
$ws.Range("C199").Value = 2026
$ws.Range("D199").Value = "MDPI AG"
$ws.Range("E199").Value = "Materials"
$ws.Range("F199").Value = "Bai, Zhaopeng; Liang, Jinsong; Ding, Chengxi; Zhou, Zimo; Luo, Man; Gu, Lin; Ma, Hong-Ping; Zhang, Qing-Chun"
$ws.Range("H199").Value = "10.3390/ma19040766"
$ws.Range("I199").Value = "https://doi.org/10.3390/ma19040766"
$ws.Range("J199").Value = "Journal"
$ws.Range("K199").Value = "n-FET"
$ws.Range("L199").Value = "Experiment"
$ws.Range("M199").Value = "Contacts"
$ws.Range("Q199").Value = "Advances and Perspectives in Gate Dielectric Thin Films for 4H-SiC MOSFETs"
$ws.Range("R199").Value = "High"
Set-DateTextCell "S199" "2026-02-16"
